$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_sheet_template")

# --- Row 1 gets slightly taller (extra wrapped line fits) ---
$ws.Rows.Item(1).RowHeight = 218.25

# --- "Allowances" header (column L, row 1): drop the ", 1/12 of annual
#     total" qualifier from the note, and give the header its own
#     (explicit black) font color, same as the other differently-colored
#     headers in the sheet. ---
$cell = $ws.Cells.Item(1, 12)

$newText = "Allowances`n(Shift, Sunday and night work plus other extra pay for difficult working conditions)"
$cell.Value2 = $newText

# Whole-cell base formatting: bold + explicit black (creates/uses the
# dedicated style for this header cell).
$cell.Font.Bold = $true
$cell.Font.Color = 0

# First run - "Allowances" (bold run)
$titleLen = "Allowances".Length
$r1 = $cell.Characters(1, $titleLen)
$r1.Font.Bold = $true
$r1.Font.Color = 0
$r1.Font.Size = 9
$r1.Font.Name = "Arial"

# Second run - the explanatory note (regular run)
$r2 = $cell.Characters($titleLen + 1, $newText.Length - $titleLen)
$r2.Font.Bold = $false
$r2.Font.Color = 0
$r2.Font.Size = 9
$r2.Font.Name = "Arial"
